# Apply cryptocurrency price/volume updates to sheet1 (ActiveSheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number but must remain
# stored as text (matching the original inlineStr cell type in the workbook).
$textForceCells = @("D5", "D6", "D7", "D10", "D12", "D16", "D20", "D21", "D24", "D25", "D27", "D29", "D33", "D35", "D36", "D38", "D40", "D41", "D43", "D45", "D47")

$updates = @(
    @{Cell = "D2"; Value = '64.575.07'},
    @{Cell = "E2"; Value = '  -1.92%  '},
    @{Cell = "D3"; Value = '3.418.08'},
    @{Cell = "E3"; Value = '  -2.20%  '},
    @{Cell = "E4"; Value = '  +0.00%  '},
    @{Cell = "D5"; Value = '571.86'},
    @{Cell = "E5"; Value = '  -1.64%  '},
    @{Cell = "D6"; Value = '157.73'},
    @{Cell = "E6"; Value = '  -2.24%  '},
    @{Cell = "D7"; Value = '0.610'},
    @{Cell = "E7"; Value = '  +0.29%  '},
    @{Cell = "E8"; Value = '  +0.07%  '},
    @{Cell = "D9"; Value = '3.418.04'},
    @{Cell = "E9"; Value = '  -2.28%  '},
    @{Cell = "D10"; Value = '7.17'},
    @{Cell = "E10"; Value = '  -2.08%  '},
    @{Cell = "E11"; Value = '  -2.84%  '},
    @{Cell = "D12"; Value = '0.438'},
    @{Cell = "E12"; Value = '  -2.39%  '},
    @{Cell = "D13"; Value = '4.007.13'},
    @{Cell = "E13"; Value = '  -2.12%  '},
    @{Cell = "E14"; Value = '  -0.32%  '},
    @{Cell = "E15"; Value = '  -4.62%  '},
    @{Cell = "D16"; Value = '27.60'},
    @{Cell = "E16"; Value = '  -4.48%  '},
    @{Cell = "D17"; Value = '64.632.19'},
    @{Cell = "E17"; Value = '  -1.78%  '},
    @{Cell = "D18"; Value = '3.423.74'},
    @{Cell = "E18"; Value = '  -2.27%  '},
    @{Cell = "E19"; Value = '  -2.49%  '},
    @{Cell = "D20"; Value = '13.78'},
    @{Cell = "E20"; Value = '  -4.08%  '},
    @{Cell = "D21"; Value = '379.55'},
    @{Cell = "E21"; Value = '  -3.28%  '},
    @{Cell = "E22"; Value = '  -3.96%  '},
    @{Cell = "E23"; Value = '  -1.49%  '},
    @{Cell = "D24"; Value = '0.999'},
    @{Cell = "E24"; Value = '  -0.32%  '},
    @{Cell = "D25"; Value = '71.95'},
    @{Cell = "E25"; Value = '  -2.21%  '},
    @{Cell = "E26"; Value = '  -5.80%  '},
    @{Cell = "D27"; Value = '9.95'},
    @{Cell = "E27"; Value = '  +1.33%  '},
    @{Cell = "D29"; Value = '1.00'},
    @{Cell = "E29"; Value = '  +0.26%  '},
    @{Cell = "E30"; Value = '  +0.50%  '},
    @{Cell = "E31"; Value = '  -3.95%  '},
    @{Cell = "E32"; Value = '  -2.82%  '},
    @{Cell = "D33"; Value = '23.17'},
    @{Cell = "E33"; Value = '  -2.63%  '},
    @{Cell = "E34"; Value = '  -1.76%  '},
    @{Cell = "D35"; Value = '1.57'},
    @{Cell = "E35"; Value = '  +1.49%  '},
    @{Cell = "D36"; Value = '160.37'},
    @{Cell = "E36"; Value = '  -1.65%  '},
    @{Cell = "E37"; Value = '  -3.19%  '},
    @{Cell = "D38"; Value = '0.0753'},
    @{Cell = "E38"; Value = '  -3.02%  '},
    @{Cell = "D39"; Value = '2.888.11'},
    @{Cell = "E39"; Value = '  -6.69%  '},
    @{Cell = "D40"; Value = '6.69'},
    @{Cell = "E40"; Value = '  +2.16%  '},
    @{Cell = "D41"; Value = '26.26'},
    @{Cell = "E41"; Value = '  -3.82%  '},
    @{Cell = "E42"; Value = '  +0.39%  '},
    @{Cell = "D43"; Value = '43.02'},
    @{Cell = "E43"; Value = '  -0.34%  '},
    @{Cell = "E44"; Value = '  -2.29%  '},
    @{Cell = "D45"; Value = '0.769'},
    @{Cell = "E45"; Value = '  -1.71%  '},
    @{Cell = "E46"; Value = '  -0.61%  '},
    @{Cell = "D47"; Value = '316.09'},
    @{Cell = "E47"; Value = '  +0.50%  '},
    @{Cell = "E48"; Value = '  -0.87%  '},
    @{Cell = "E49"; Value = '  -5.78%  '},
    @{Cell = "E50"; Value = '  -2.28%  '},
    @{Cell = "E51"; Value = '  -3.39%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($textForceCells -contains $u.Cell) {
        # Preserve text storage: temporarily force a text number format so
        # Excel does not reinterpret the numeric-looking string as a number,
        # then restore the cell's original style/format.
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = $origStyle
    } else {
        $rng.Value = $u.Value
    }
}
